# Update crypto price/volume figures per the Jan 20 2024 09:06:51 UTC data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.38"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "2.475.17"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'313.63"

$ws.Range("D6").Value = "'92.01"
$ws.Range("E6").Value = "  -2.95%  "

$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10").Value = "'32.58"
$ws.Range("E10").Value = "  -2.99%  "

$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").Value = "2.855.87"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "'16.38"
$ws.Range("E14").Value = "  +8.25%  "

$ws.Range("D15").Value = "'6.89"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").Value = "2.408.32"
$ws.Range("E16").Value = "  -2.72%  "

$ws.Range("D17").Value = "'0.771"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("D18").Value = "41.534.79"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("D21").Value = "'71.84"
$ws.Range("E21").Value = "  +4.74%  "

$ws.Range("D22").Value = "'11.12"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").Value = "'236.39"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "'24.89"
$ws.Range("E27").Value = "  +3.42%  "

$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").Value = "'35.65"
$ws.Range("E30").Value = "  -2.49%  "

$ws.Range("D31").Value = "'156.19"
$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("D32").Value = "'5.47"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").Value = "'0.0760"
$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  -8.07%  "

$ws.Range("D37").Value = "'2.90"
$ws.Range("E37").Value = "  -6.19%  "

$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("E41").Value = "  -5.63%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").Value = "1.960.42"
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("D44").Value = "'0.0285"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("D45").Value = "'18.67"
$ws.Range("E45").Value = "  -5.08%  "

$ws.Range("D46").Value = "'2.94"
$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").Value = "'9.06"
$ws.Range("E47").Value = "  +4.15%  "

$ws.Range("D48").Value = "2.713.74"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").Value = "'97.58"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'67.56"
$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("D51").Value = "'72.21"
$ws.Range("E51").Value = "  -3.06%  "
